$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number that is bumped by one
# day for rows 2-14 as part of the automatic daily update.
for ($row = 2; $row -le 14; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46074) {
        $cell.Value2 = 46075
    }
}
